$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2.24
$ws.Range("G2").Value = 2.46
$ws.Range("H2").Value = 3.85
$ws.Range("I2").Value = 4.5
$ws.Range("J2").Value = 2.86
$ws.Range("K2").Value = 3.2
$ws.Range("L2").Value = 1.66
$ws.Range("M2").Value = 1.15
$ws.Range("N2").Value = 2.34
$ws.Range("O2").Value = 1.65
$ws.Range("P2").Value = 1.42
$ws.Range("Q2").Value = 2.98
$ws.Range("R2").Value = 1.14
$ws.Range("S2").Value = 6.4
$ws.Range("T2").Value = 2.34
$ws.Range("U2").Value = 1.61
$ws.Range("V2").Value = 1.29
$ws.Range("W2").Value = 1.69
$ws.Range("Y2").Value = 23
$ws.Range("AB2").Value = 12
$ws.Range("AC2").Value = 7.6
$ws.Range("AF2").Value = 24
$ws.Range("AG2").Value = 24
$ws.Range("AH2").Value = 990
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 95
$ws.Range("AL2").Value = 1000
$ws.Range("F3").Value = 2.66
$ws.Range("G3").Value = 3
$ws.Range("H3").Value = 3.1
$ws.Range("I3").Value = 3.5
$ws.Range("L3").Value = 1.69
$ws.Range("N3").Value = 2.28
$ws.Range("O3").Value = 1.66
$ws.Range("P3").Value = 1.41
$ws.Range("Q3").Value = 3
$ws.Range("S3").Value = 6.6
$ws.Range("T3").Value = 2.28
$ws.Range("U3").Value = 1.64
$ws.Range("V3").Value = 1.4
$ws.Range("W3").Value = 1.51
$ws.Range("X3").Value = 12.5
$ws.Range("Y3").Value = 14
$ws.Range("Z3").Value = 48
$ws.Range("AB3").Value = 15
$ws.Range("AE3").Value = 300
$ws.Range("AF3").Value = 65
$ws.Range("AG3").Value = 30
$ws.Range("AH3").Value = 990
$ws.Range("AJ3").Value = 220
$ws.Range("F4").Value = 1.15
$ws.Range("G4").Value = 1.18
$ws.Range("H4").Value = 19
$ws.Range("J4").Value = 10
$ws.Range("K4").Value = 12
$ws.Range("L4").Value = 1.18
$ws.Range("N4").Value = 9.199999999999999
$ws.Range("O4").Value = 1.09
$ws.Range("P4").Value = 3.95
$ws.Range("Q4").Value = 1.28
$ws.Range("R4").Value = 2.18
$ws.Range("S4").Value = 1.75
$ws.Range("T4").Value = 1.96
$ws.Range("W4").Value = 6.8
$ws.Range("X4").Value = 75
$ws.Range("Z4").Value = 290
$ws.Range("AB4").Value = 17
$ws.Range("AC4").Value = 30
$ws.Range("AD4").Value = 990
$ws.Range("AE4").Value = 380
$ws.Range("AF4").Value = 11
$ws.Range("AH4").Value = 44
$ws.Range("AN4").Value = 2.64
$ws.Range("F5").Value = 1.92
$ws.Range("G5").Value = 2.08
$ws.Range("H5").Value = 5.2
$ws.Range("J5").Value = 2.9
$ws.Range("K5").Value = 3.3
$ws.Range("L5").Value = 1.67
$ws.Range("M5").Value = 1.16
$ws.Range("N5").Value = 2.26
$ws.Range("O5").Value = 1.69
$ws.Range("P5").Value = 1.4
$ws.Range("Q5").Value = 3.1
$ws.Range("R5").Value = 1.14
$ws.Range("S5").Value = 7
$ws.Range("T5").Value = 2.52
$ws.Range("U5").Value = 1.53
$ws.Range("W5").Value = 1.92
$ws.Range("X5").Value = 7.2
$ws.Range("Y5").Value = 12.5
$ws.Range("AB5").Value = 5.7
$ws.Range("AD5").Value = 27
$ws.Range("AE5").Value = 1000
$ws.Range("AG5").Value = 13
$ws.Range("AH5").Value = 34
$ws.Range("AJ5").Value = 27
$ws.Range("AK5").Value = 34
$ws.Range("AL5").Value = 100
$ws.Range("AN5").Value = 34
$ws.Range("L6").Value = 1.58
$ws.Range("N6").Value = 2.84
$ws.Range("O6").Value = 1.53
$ws.Range("Q6").Value = 2.6
$ws.Range("R6").Value = 1.21
$ws.Range("T6").Value = 2.1
$ws.Range("U6").Value = 1.87
$ws.Range("X6").Value = 9
$ws.Range("Y6").Value = 9.4
$ws.Range("AB6").Value = 7.8
$ws.Range("AD6").Value = 15
$ws.Range("AN6").Value = 38
$ws.Range("F7").Value = 1.57
$ws.Range("G7").Value = 1.61
$ws.Range("I7").Value = 5.8
$ws.Range("N7").Value = 7.4
$ws.Range("P7").Value = 3.2
$ws.Range("Q7").Value = 1.4
$ws.Range("R7").Value = 1.91
$ws.Range("S7").Value = 2.04
$ws.Range("U7").Value = 2.62
$ws.Range("W7").Value = 2.62
$ws.Range("X7").Value = 44
$ws.Range("Y7").Value = 95
$ws.Range("AB7").Value = 16
$ws.Range("AC7").Value = 13
$ws.Range("AE7").Value = 60
$ws.Range("AF7").Value = 14.5
$ws.Range("AG7").Value = 10.5
$ws.Range("AH7").Value = 17.5
$ws.Range("AI7").Value = 48
$ws.Range("AJ7").Value = 17.5
$ws.Range("AL7").Value = 23
$ws.Range("AN7").Value = 5.2
$ws.Range("AO7").Value = 600
$ws.Range("F8").Value = 4.8
$ws.Range("G8").Value = 5.3
$ws.Range("H8").Value = 1.83
$ws.Range("I8").Value = 1.91
$ws.Range("P8").Value = 1.93
$ws.Range("R8").Value = 1.36
$ws.Range("V8").Value = 2.08
$ws.Range("X8").Value = 15
$ws.Range("Z8").Value = 11.5
$ws.Range("AA8").Value = 38
$ws.Range("AI8").Value = 980
$ws.Range("AO8").Value = 13.5
$ws.Range("F9").Value = 1.69
$ws.Range("H9").Value = 6.2
$ws.Range("I9").Value = 7.6
$ws.Range("K9").Value = 3.8
$ws.Range("L9").Value = 1.57
$ws.Range("N9").Value = 2.62
$ws.Range("O9").Value = 1.52
$ws.Range("Q9").Value = 2.54
$ws.Range("R9").Value = 1.19
$ws.Range("S9").Value = 5.3
$ws.Range("T9").Value = 2.32
$ws.Range("V9").Value = 1.15
$ws.Range("X9").Value = 9.4
$ws.Range("Y9").Value = 17
$ws.Range("AC9").Value = 8.6
$ws.Range("AD9").Value = 30
$ws.Range("AF9").Value = 9.199999999999999
$ws.Range("AJ9").Value = 19.5
$ws.Range("AK9").Value = 26
$ws.Range("AN9").Value = 20
$ws.Range("H10").Value = 9.800000000000001
$ws.Range("I10").Value = 10
$ws.Range("O10").Value = 1.41
$ws.Range("Q10").Value = 2.16
$ws.Range("R10").Value = 1.29
$ws.Range("S10").Value = 4.1
$ws.Range("T10").Value = 2.46
$ws.Range("Y10").Value = 25
$ws.Range("AA10").Value = 420
$ws.Range("AC10").Value = 10.5
$ws.Range("AE10").Value = 210
$ws.Range("AF10").Value = 7.2
$ws.Range("AO10").Value = 360
$ws.Range("F11").Value = 1.62
$ws.Range("G11").Value = 1.63
$ws.Range("H11").Value = 6.2
$ws.Range("I11").Value = 6.4
$ws.Range("J11").Value = 4.4
$ws.Range("K11").Value = 4.5
$ws.Range("L11").Value = 1.41
$ws.Range("N11").Value = 4.1
$ws.Range("O11").Value = 1.3
$ws.Range("Q11").Value = 1.92
$ws.Range("R11").Value = 1.42
$ws.Range("S11").Value = 3.3
$ws.Range("T11").Value = 1.93
$ws.Range("U11").Value = 2.06
$ws.Range("V11").Value = 1.18
$ws.Range("W11").Value = 2.58
$ws.Range("X11").Value = 17
$ws.Range("Y11").Value = 20
$ws.Range("Z11").Value = 50
$ws.Range("AA11").Value = 180
$ws.Range("AB11").Value = 8
$ws.Range("AD11").Value = 23
$ws.Range("AK11").Value = 16
$ws.Range("F12").Value = 2.12
$ws.Range("G12").Value = 2.24
$ws.Range("H12").Value = 3.85
$ws.Range("I12").Value = 4.3
$ws.Range("J12").Value = 3.2
$ws.Range("K12").Value = 3.5
$ws.Range("L12").Value = 1.53
$ws.Range("Q12").Value = 2.32
$ws.Range("S12").Value = 4.4
$ws.Range("T12").Value = 1.97
$ws.Range("U12").Value = 1.88
$ws.Range("W12").Value = 1.81
$ws.Range("AB12").Value = 8.199999999999999
$ws.Range("AE12").Value = 240
$ws.Range("AI12").Value = 75
$ws.Range("AJ12").Value = 30
$ws.Range("AN12").Value = 30
